# Update the crypto price/volume table (columns D "Price" and E "Volume(1h)")
# with the latest scraped values. Cells that look like plain numbers (e.g.
# "214.79") are written with a leading apostrophe so Excel stores them as
# text - matching the original workbook, which keeps every D/E cell as a
# string (some prices use "." as a thousands separator, e.g. "26.951.25").
# The Style reset afterwards clears the resulting quote-prefix formatting
# so the cell keeps its original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.951.25'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').Value = '1.650.87'
$ws.Range('E3').Value = '  +2.84%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'214.79"
$ws.Range('E5').Value = '  +1.32%  '
$ws.Range('D6').Value = "'0.511"
$ws.Range('E6').Value = '  +2.44%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +2.65%  '
$ws.Range('E9').Value = '  +1.64%  '
$ws.Range('E10').Value = '  +4.82%  '
$ws.Range('D11').Value = "'0.0879"
$ws.Range('E11').Value = '  +2.88%  '
$ws.Range('D12').Value = '1.884.90'
$ws.Range('E12').Value = '  +2.94%  '
$ws.Range('D13').Value = '1.641.88'
$ws.Range('E13').Value = '  +2.32%  '
$ws.Range('E15').Value = '  +2.68%  '
$ws.Range('D16').Value = "'65.06"
$ws.Range('E16').Value = '  +2.81%  '
$ws.Range('D17').Value = '26.946.27'
$ws.Range('E17').Value = '  +2.17%  '
$ws.Range('D18').Value = "'235.65"
$ws.Range('E18').Value = '  +2.09%  '
$ws.Range('D19').Value = '0.0₃0733'
$ws.Range('E19').Value = '  +1.24%  '
$ws.Range('D20').Value = "'7.74"
$ws.Range('E20').Value = '  +1.10%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  +3.11%  '
$ws.Range('D23').Value = "'9.32"
$ws.Range('E23').Value = '  +4.02%  '
$ws.Range('E24').Value = '  +2.11%  '
$ws.Range('D25').Value = "'145.36"
$ws.Range('E25').Value = '  -1.07%  '
$ws.Range('E26').Value = '  +2.15%  '
$ws.Range('E27').Value = '  +0.82%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  +2.65%  '
$ws.Range('D30').Value = "'0.0496"
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('E31').Value = '  +1.92%  '
$ws.Range('D32').Value = '1.557.79'
$ws.Range('E32').Value = '  +4.47%  '
$ws.Range('E33').Value = '  +2.71%  '
$ws.Range('D34').Value = "'3.08"
$ws.Range('E34').Value = '  +4.82%  '
$ws.Range('E35').Value = '  +9.43%  '
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D37').Value = "'0.585"
$ws.Range('E37').Value = '  +4.19%  '
$ws.Range('E38').Value = '  +8.97%  '
$ws.Range('E39').Value = '  +2.90%  '
$ws.Range('E40').Value = '  +3.28%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').Value = "'65.92"
$ws.Range('E42').Value = '  +8.36%  '
$ws.Range('E43').Value = '  +2.24%  '
$ws.Range('D44').Value = '1.791.73'
$ws.Range('E45').Value = '  +2.25%  '
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('D47').Value = "'89.96"
$ws.Range('E47').Value = '  +0.64%  '
$ws.Range('E48').Value = '  +1.80%  '
$ws.Range('D49').Value = "'0.0989"
$ws.Range('E49').Value = '  +3.03%  '
$ws.Range('D51').Value = "'7.64"
$ws.Range('E51').Value = '  +3.02%  '

$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
